# Refresh crypto price/volume data (and re-sort a few coin rows) per upstream feed update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells below keep a trailing zero / fixed-point format that plain numeric entry would drop,
# so mark them as Text before writing the value (same as typing into a Text-formatted cell).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'

$ws.Range('D2').Value = '29.948.97'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.899.82'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('D4').Value = '0.9981'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '0.7900'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').Value = '244.33'
$ws.Range('E6').Value = '  +1.11%  '
$ws.Range('D7').Value = '0.9983'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').Value = '0.3165'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '25.61'
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').Value = '0.07362'
$ws.Range('E10').Value = '  +4.68%  '
$ws.Range('D11').Value = '0.08126'
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('D12').Value = '0.7762'
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('D13').Value = '5.515'
$ws.Range('E13').Value = '  +3.87%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').Value = '94.08'
$ws.Range('E14').Value = '  +1.82%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.791.87'
$ws.Range('E15').Value = '  -6.73%  '
$ws.Range('D16').Value = '6.249'
$ws.Range('E16').Value = '  +5.33%  '
$ws.Range('D17').Value = '29.884.15'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('D18').Value = '14.01'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('D19').Value = '246.65'
$ws.Range('D20').Value = '0.000007867'
$ws.Range('E20').Value = '  +2.21%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').Value = '8.170'
$ws.Range('E21').Value = '  +0.57%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '0.9982'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('B23').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D23').Value = '2.095.00'
$ws.Range('E23').Value = '  -2.35%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Value = '0.9980'
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('D25').Value = '0.1605'
$ws.Range('E25').Value = '  -2.54%  '
$ws.Range('D26').Value = '9.503'
$ws.Range('E26').Value = '  +2.10%  '
$ws.Range('D27').Value = '162.54'
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('D28').Value = '18.84'
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('D29').Value = '2.044'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').Value = '1.458'
$ws.Range('E30').Value = '  +5.46%  '
$ws.Range('D31').Value = '1.548'
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('D32').Value = '4.509'
$ws.Range('E32').Value = '  +1.94%  '
$ws.Range('D33').Value = '0.05623'
$ws.Range('E33').Value = '  -1.60%  '
$ws.Range('D34').Value = '4.102'
$ws.Range('E34').Value = '  +0.57%  '
$ws.Range('D35').Value = '1.257'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').Value = '0.7586'
$ws.Range('E36').Value = '  +2.98%  '
$ws.Range('E38').Value = '  +1.44%  '
$ws.Range('D39').Value = '0.01943'
$ws.Range('E39').Value = '  +1.80%  '
$ws.Range('D40').Value = '2.798'
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('D41').Value = '1.150.10'
$ws.Range('E41').Value = '  +12.06%  '
$ws.Range('D42').Value = '0.4480'
$ws.Range('E42').Value = '  +1.78%  '
$ws.Range('D43').Value = '74.03'
$ws.Range('E43').Value = '  +2.07%  '
$ws.Range('D44').Value = '5.975'
$ws.Range('E44').Value = '  +2.48%  '
$ws.Range('D45').Value = '0.8578'
$ws.Range('E45').Value = '  +2.20%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '1.906'
$ws.Range('E46').Value = '  +3.35%  '
$ws.Range('B47').Value = 'SynthetixNetwork'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D47').Value = '3.170'
$ws.Range('E47').Value = '  +8.62%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').Value = '0.9981'
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('D49').Value = '102.22'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').Value = '7.551'
$ws.Range('E50').Value = '  +1.41%  '
$ws.Range('D51').Value = '9.804'
$ws.Range('E51').Value = '  -0.99%  '
